# "Placed LEDs with 1:60 scale"
#
# The sheet drives an LED placement grid from measured X/Y positions (columns
# B and C) combined with a Scale (M1), X/Y correction factors (M2/M3) and
# X/Y offsets (M4/M5). This edit:
#   1. Updates the Scale value (M1) from 50 to 80.
#   2. Re-derives the scaled output formulas (I and J columns) so the offset
#      is applied AFTER the correction-factor scaling, instead of before:
#        I: (F+offset)*corr  ->  (F*corr)+offset
#        J: (G+offset)*corr  ->  (G*corr)+offset

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Scale
$ws.Range("M1").Value = 80

# 2. Re-derive formulas, keeping the existing shared-formula groupings
#    (row 1 stand-alone, rows 2-65, rows 66-82) so the edit mirrors the
#    original authoring structure as closely as possible.
$ws.Range("I1").Formula = "=(F1*`$M`$2)+`$M`$4"
$ws.Range("I2:I65").Formula = "=(F2*`$M`$2)+`$M`$4"
$ws.Range("I66:I82").Formula = "=(F66*`$M`$2)+`$M`$4"

$ws.Range("J1").Formula = "=(G1*`$M`$3)+`$M`$5"
$ws.Range("J2:J65").Formula = "=(G2*`$M`$3)+`$M`$5"
$ws.Range("J66:J82").Formula = "=(G66*`$M`$3)+`$M`$5"

# Re-entering the formulas stamps a stray numeric style onto the touched
# cells; restore them to the original (unstyled) Normal style.
$ws.Range("I1:I82").Style = "Normal"
$ws.Range("J1:J82").Style = "Normal"
